$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (STRASSE, HAUSNR) before the old column J (PLZ),
# shifting PLZ/WOHNORT/VERMÖGEN/HASEL/HASSH two columns to the right.
$ws.Range("J1:K1").EntireColumn.Insert()

# Match the column width used by the rest of the "address" columns.
$ws.Range("J1:K1").ColumnWidth = 22.33

# Header row (set HAUSNR before STRASSE so the shared-string table order
# matches the authored workbook)
$ws.Range("K1").Value = "HAUSNR"
$ws.Range("J1").Value = "STRASSE"

# Row 2 - Aaron Ackermann
$ws.Range("J2").Value = "Ackerstrasse"
$ws.Range("K2").Value = 11

# Row 3 - Berta Brunner
$ws.Range("J3").Value = "Bertastrasse"
$ws.Range("K3").Value = 22

# Row 4 - Claudio Christen
$ws.Range("J4").Value = "Clausiensteig"
$ws.Range("K4").Value = "3c"

# Row 6 - Dario Dachs
$ws.Range("J6").Value = "Dammweg"
$ws.Range("K6").Value = 4

# Row 7 - Elsa Eris
$ws.Range("J7").Value = "Erismannstrasse"
$ws.Range("K7").Value = 505

# Row 8 - Fiona Fichter
$ws.Range("J8").Value = "Floragasse"
$ws.Range("K8").Value = "6 f"

# Update the saved selection to match the authored workbook state.
$ws.Range("K12").Select() | Out-Null
